$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")

# Clear the expense-row details that were removed from the sheet
# (values only; cell formatting/styles are left untouched, matching
# a plain Delete-key / ClearContents edit over these ranges).
$ws.Range("B8:F9").ClearContents()
$ws.Range("B11:F16").ClearContents()
$ws.Range("B24:F25").ClearContents()

# Update the current selection to match the cleared block.
$ws.Range("B24:F25").Select()
